# Apply updated cryptocurrency price/volume values to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.605.68"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.594.89"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'211.10"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.24%  "
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").Value = "'19.36"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "'0.0837"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "1.619.53"
$ws.Range("E13").Value = "  +2.26%  "
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "'0.521"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "'64.53"
$ws.Range("E16").Value = "  -0.32%  "
$ws.Range("D17").Value = "26.588.36"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").Value = "'208.50"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").Value = "'6.95"
$ws.Range("E21").Value = "  +3.49%  "
$ws.Range("D22").Value = "'4.25"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("E23").Value = "  -1.94%  "
$ws.Range("D24").Value = "'8.86"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'145.07"
$ws.Range("E25").Value = "  -1.10%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'7.09"
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "'0.0504"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("E32").Value = "  +0.41%  "
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D35").Value = "1.278.37"
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("D36").Value = "'2.45"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").Value = "'0.841"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("D42").Value = "'2.19"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("D43").Value = "'0.785"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("D44").Value = "'63.99"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "1.730.00"
$ws.Range("D46").Value = "'0.913"
$ws.Range("E46").Value = "  +9.09%  "
$ws.Range("D47").Value = "'89.54"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  -0.99%  "
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("E50").Value = "  +4.00%  "
$ws.Range("E51").Value = "  +0.40%  "
